# Apply cryptos list update (prices + 1h volumes), per commit
# "Updated cryptos list on Sun May 14 08:36:52 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.777.95"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "1.864.83"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").Value = "  +2.69%  "
$ws.Range("D5").Value = "'324.97"
$ws.Range("E5").Value = "  +3.78%  "
$ws.Range("D6").Value = "'1.034"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("D7").Value = "'0.4415"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("D8").Value = "'0.3795"
$ws.Range("E8").Value = "  +2.66%  "
$ws.Range("D9").Value = "'0.07462"
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").Value = "'0.8846"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "'21.79"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").Value = "1.885.06"
$ws.Range("E12").Value = "  -11.99%  "
$ws.Range("D13").Value = "'5.557"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "'6.753"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "'0.07229"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("D16").Value = "'83.72"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "'1.040"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").Value = "'0.000009139"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").Value = "'15.55"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").Value = "27.789.14"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "'5.312"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").Value = "'11.40"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").Value = "'1.968"
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("D25").Value = "'158.33"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").Value = "'18.87"
$ws.Range("E26").Value = "  +2.71%  "
$ws.Range("D27").Value = "'1.990"
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("D28").Value = "'5.322"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("D29").Value = "'117.66"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "'0.09107"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("B31").Value = "ARBITRUM"
$ws.Range("C31").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D31").Value = "'1.218"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.7756"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("D33").Value = "'3.044"
$ws.Range("E33").Value = "  +8.01%  "
$ws.Range("D34").Value = "'4.588"
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("D35").Value = "'1.036"
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").Value = "'1.167"
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("D37").Value = "'0.01991"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").Value = "'0.05351"
$ws.Range("E38").Value = "  +2.09%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.844"
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5198"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").Value = "'0.1695"
$ws.Range("E41").Value = "  +2.49%  "
$ws.Range("D42").Value = "'6.895"
$ws.Range("E42").Value = "  +5.81%  "
$ws.Range("D43").Value = "'8.727"
$ws.Range("E43").Value = "  +4.97%  "
$ws.Range("D44").Value = "'109.64"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").Value = "'10.71"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").Value = "'1.722"
$ws.Range("E46").Value = "  +4.25%  "
$ws.Range("D47").Value = "'0.4709"
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("D48").Value = "'0.06435"
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("D49").Value = "'1.884"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("D50").Value = "'39.92"
$ws.Range("E50").Value = "  +4.72%  "
$ws.Range("D51").Value = "'64.55"
$ws.Range("E51").Value = "  +1.17%  "
